# Add a "Natural Earth" country-name mapping column (B) to the Country
# lookup sheet, mapping each EMPD country name (column A) onto the
# corresponding Natural Earth country name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Country")

# Row -> Natural Earth country name. Ordered to match the row order of the
# sheet (top to bottom); rows 81 (Black Sea), 82 (Dead Sea) and 84 (Adriatic
# Sea) are not real countries and are intentionally left without a mapping.
$map = [ordered]@{
    2  = "Albania"
    3  = "Algeria"
    4  = "Andorra"
    5  = "Armenia"
    6  = "Austria"
    7  = "Azerbaijan"
    8  = "Bahrain"
    9  = "Belarus"
    10 = "Belgium"
    11 = "Bosnia and Herzegovina"
    12 = "Bulgaria"
    13 = "Cape Verde"
    14 = "Croatia"
    15 = "Cyprus"
    16 = "Czech Republic"
    17 = "Denmark"
    18 = "Egypt"
    19 = "Eritrea"
    20 = "Estonia"
    21 = "Ethiopia"
    22 = "Faroe Islands"
    23 = "Finland"
    24 = "France"
    25 = "Georgia"
    26 = "Germany"
    27 = "Gibraltar"
    28 = "Greece"
    29 = "Greenland"
    30 = "Guernsey"
    31 = "Hungary"
    32 = "Iceland"
    33 = "Iran"
    34 = "Iraq"
    35 = "Ireland"
    36 = "Isle of Man"
    37 = "Israel"
    38 = "Italy"
    39 = "India"
    40 = "Japan"
    41 = "Jersey"
    42 = "Kazakhstan"
    43 = "Kuwait"
    44 = "Kyrgyzstan"
    45 = "Latvia"
    46 = "Lebanon"
    47 = "Libya"
    48 = "Liechtenstein"
    49 = "Lithuania"
    50 = "Luxembourg"
    51 = "Macedonia"
    52 = "Malta"
    53 = "Moldova"
    54 = "Monaco"
    55 = "Morocco"
    56 = "Netherlands"
    57 = "Norway"
    58 = "Oman"
    59 = "Palestine"
    60 = "Poland"
    61 = "Portugal"
    62 = "Qatar"
    63 = "Romania"
    64 = "Russia"
    65 = "San Marino"
    66 = "Saudi Arabia"
    67 = "Montenegro"
    68 = "Slovakia"
    69 = "Slovenia"
    70 = "Spain"
    71 = "Sweden"
    72 = "Switzerland"
    73 = "Syria"
    74 = "Tunisia"
    75 = "Turkey"
    76 = "Ukraine"
    77 = "United Arab Emirates"
    78 = "United Kingdom"
    79 = "Yemen"
    80 = "Jordan"
    83 = "Turkmenistan"
    85 = "Norway"
    86 = "China"
}

foreach ($row in $map.Keys) {
    $ws.Cells.Item([int]$row, 2).Value = $map[$row]
}

# Header for the new column, added last (so it becomes the final new shared
# string), matching the bold/quote-prefixed header style already used by A1.
$ws.Cells.Item(1, 2).Value = "Natural Earth"
$ws.Range("A1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection/tab to the Country sheet, cell B1 - this also
# updates the workbook's active-tab bookkeeping and drops the previous
# scroll position / selection on this sheet.
$ws.Range("B1").Select()
